# "Logica de interaccion con el select"
# Update the VACANCY value selected for the applicant, widen the VACANCY
# column so the longer text fits, and move the active selection to reflect
# the user interacting with the (vacancy) select/dropdown on the
# Recruitment sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recruitment")
$ws.Activate()

# The applicant's selected vacancy changes from "QA LEAD" to
# "Junior Account Assistant".
$ws.Range("E2").Value = "Junior Account Assistant"

# Widen column E (VACANCY) so the new, longer value is fully visible.
$ws.Columns.Item(5).ColumnWidth = 19.65

# Reflect the user's click/selection interacting with the select, moving
# the active cell from K6 to E5.
$ws.Range("E5").Select()
